# Apply the "new version with timestamp" edit:
#   - Insert a new data row (item #8 "سرنجات 3 سم") after the current last
#     data row (item #7), re-using that row's layout/styles.
#   - Update item #7's H/L/N cells with its new values.
#   - Recompute the totals-row sum to include the new row.
# All subsequent rows (grand-total row, footer row) shift down by one and
# keep their original formatting/content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for the sheet's row 11 (pushes old row 11 -> 12, 12 -> 13) ---
$ws.Rows.Item(11).Insert()

# --- Duplicate row 10's formatting (fonts/fills/borders/number formats) onto
#     the freshly inserted row 11, without touching its (still empty) values ---
$ws.Range("A10:N10").Copy()
$ws.Range("A11:N11").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row heights: new row 11 matches other data rows; total row 12 grows slightly
$ws.Rows.Item(11).RowHeight = 24.75
$ws.Rows.Item(12).RowHeight = 26.25

# Re-create the merges for the new row 11 (mirrors B10:G10 / H10:K10 / L10:M10)
$ws.Range("B11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()

# --- Update existing row 10 (item #7) values ---
$ws.Range("H10").Value = "0:0"
$ws.Range("L10").Value = 30
$ws.Range("N10").Value = "1:0"

# --- Populate the new row 11 (item #8) ---
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "سرنجات 3 سم"
$ws.Range("H11").Value = "-10:0"
$ws.Range("L11").Value = 20
$ws.Range("N11").Value = "10:0"

# --- Update the grand-total row (now row 12) to include the new row's total ---
$ws.Range("K12").Value = 294.36
